$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "Price" (column D) values are stored as text in the source data (some
# contain thousand-separator dots like "63.225.08"), so force text format on
# column D before writing, to stop Excel auto-converting plain numeric-looking
# strings (e.g. "327.00" -> 327) and dropping significant trailing zeros.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.225.08"
$ws.Range("E2").Value = "  +1.79%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.450.46"
$ws.Range("E3").Value = "  +1.06%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.66"
$ws.Range("E5").Value = "  +1.56%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.27"
$ws.Range("E6").Value = "  +2.25%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.536"
$ws.Range("E8").Value = "  +0.90%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.453.03"
$ws.Range("E9").Value = "  +1.20%  "

# Row 10
$ws.Range("E10").Value = "  +2.58%  "

# Row 11
$ws.Range("E11").Value = "  +0.22%  "

# Row 12
$ws.Range("E12").Value = "  +2.28%  "

# Row 13
$ws.Range("E13").Value = "  +1.03%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.02"
$ws.Range("E14").Value = "  +3.10%  "

# Row 15
$ws.Range("E15").Value = "  +3.82%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.879.92"
$ws.Range("E16").Value = "  +0.72%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.080.66"
$ws.Range("E17").Value = "  +1.78%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.452.26"
$ws.Range("E18").Value = "  +0.75%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("E19").Value = "  +0.53%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.26"
$ws.Range("E20").Value = "  +6.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.00"
$ws.Range("E21").Value = "  +0.87%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").Value = "  +1.41%  "

# Row 23
$ws.Range("E23").Value = "  +11.69%  "

# Row 24
$ws.Range("E24").Value = "  -0.16%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.78"
$ws.Range("E25").Value = "  -0.79%  "

# Row 26
$ws.Range("E26").Value = "  +11.75%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.86"
$ws.Range("E27").Value = "  +0.45%  "

# Row 28
$ws.Range("E28").Value = "  +10.78%  "

# Row 29
$ws.Range("E29").Value = "  +1.94%  "

# Row 30
$ws.Range("E30").Value = "  +7.19%  "

# Row 31
$ws.Range("E31").Value = "  -0.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.31"
$ws.Range("E32").Value = "  +1.04%  "

# Row 33
$ws.Range("E33").Value = "  -1.82%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.90"
$ws.Range("E34").Value = "  +1.87%  "

# Row 35
$ws.Range("E35").Value = "  +7.57%  "

# Row 36
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("E38").Value = "  +0.74%  "

# Row 39
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.78"
$ws.Range("E39").Value = "  +0.45%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.40"
$ws.Range("E40").Value = "  -1.48%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.82"
$ws.Range("E41").Value = "  +0.23%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "145.87"
$ws.Range("E42").Value = "  -4.23%  "

# Row 43
$ws.Range("E43").Value = "  +16.74%  "

# Row 44
$ws.Range("E44").Value = "  -0.11%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.62"
$ws.Range("E45").Value = "  +0.55%  "

# Row 46
$ws.Range("E46").Value = "  +2.81%  "

# Row 47
$ws.Range("E47").Value = "  +2.01%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.82"
$ws.Range("E48").Value = "  +4.30%  "

# Row 49
$ws.Range("E49").Value = "  +0.93%  "

# Row 50
$ws.Range("E50").Value = "  +3.50%  "

# Row 51
$ws.Range("E51").Value = "  +0.44%  "
